$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "840×9=" "377×9="
Replace-Text "832×9=" "547×4="
Replace-Text "223×4=" "391×8="
Replace-Text "858×9=" "517×3="
Replace-Text "105×9=" "794×7="
Replace-Text "523×6=" "835×4="
Replace-Text "528×9=" "609×3="
Replace-Text "860×5=" "428×3="
Replace-Text "921×3=" "981×5="
Replace-Text "485×6=" "137×6="
Replace-Text "383×6=" "780×4="
Replace-Text "144×4=" "401×5="
Replace-Text "419×6=" "420×8="
Replace-Text "164×5=" "371×2="
Replace-Text "401×4=" "184×2="
Replace-Text "975×4=" "359×3="
Replace-Text "248×2=" "687×2="
Replace-Text "491×4=" "647×8="
Replace-Text "125×3=" "143×7="
Replace-Text "703×6=" "689×6="
Replace-Text "242×3=" "505×6="
Replace-Text "536×7=" "550×4="
Replace-Text "456×4=" "872×7="
Replace-Text "305×3=" "704×8="
Replace-Text "149×3=" "797×3="

Write-Output "Replacements complete"
